$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.471.22"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "2.285.88"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'310.15"
$ws.Range("E5").Value = "  -4.11%  "

# Row 6
$ws.Range("D6").Value = "'102.98"
$ws.Range("E6").Value = "  -1.38%  "

# Row 7
$ws.Range("E7").Value = "  -1.22%  "

# Row 8
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("E9").Value = "  -1.57%  "

# Row 10
$ws.Range("D10").Value = "'38.64"
$ws.Range("E10").Value = "  -3.88%  "

# Row 11
$ws.Range("E11").Value = "  -1.11%  "

# Row 12
$ws.Range("D12").Value = "'8.19"
$ws.Range("E12").Value = "  -4.02%  "

# Row 13
$ws.Range("E13").Value = "  +0.57%  "

# Row 14
$ws.Range("D14").Value = "'0.967"
$ws.Range("E14").Value = "  -0.46%  "

# Row 15
$ws.Range("D15").Value = "'15.23"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("D16").Value = "2.632.18"
$ws.Range("E16").Value = "  -0.53%  "

# Row 17
$ws.Range("D17").Value = "2.284.36"
$ws.Range("E17").Value = "  -0.22%  "

# Row 18
$ws.Range("D18").Value = "42.697.96"
$ws.Range("E18").Value = "  +0.19%  "

# Row 19
$ws.Range("D19").Value = "'7.27"
$ws.Range("E19").Value = "  -2.45%  "

# Row 20
$ws.Range("E20").Value = "  -1.50%  "

# Row 21
$ws.Range("D21").Value = "'13.40"
$ws.Range("E21").Value = "  -0.49%  "

# Row 22
$ws.Range("D22").Value = "'73.14"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").Value = "'268.40"
$ws.Range("E23").Value = "  -0.58%  "

# Row 24
$ws.Range("D24").Value = "'3.39"
$ws.Range("E24").Value = "  -5.68%  "

# Row 25
$ws.Range("D25").Value = "'2.16"
$ws.Range("E25").Value = "  -3.19%  "

# Row 26
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("D27").Value = "'10.72"
$ws.Range("E27").Value = "  -1.80%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  -0.76%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'6.91"
$ws.Range("E29").Value = "  +11.96%  "

# Row 30
$ws.Range("D30").Value = "'22.27"
$ws.Range("E30").Value = "  -1.29%  "

# Row 31
$ws.Range("D31").Value = "'35.77"
$ws.Range("E31").Value = "  -7.31%  "

# Row 32
$ws.Range("D32").Value = "'164.29"
$ws.Range("E32").Value = "  -0.79%  "

# Row 33
$ws.Range("D33").Value = "'0.0844"
$ws.Range("E33").Value = "  -4.26%  "

# Row 34
$ws.Range("E34").Value = "  -2.82%  "

# Row 35
$ws.Range("D35").Value = "'2.54"
$ws.Range("E35").Value = "  +0.67%  "

# Row 36
$ws.Range("E36").Value = "  -2.98%  "

# Row 37
$ws.Range("E37").Value = "  -3.19%  "

# Row 38
$ws.Range("E38").Value = "  -3.37%  "

# Row 39
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("D40").Value = "'3.60"
$ws.Range("E40").Value = "  -3.93%  "

# Row 41
$ws.Range("D41").Value = "'111.29"
$ws.Range("E41").Value = "  +19.44%  "

# Row 42
$ws.Range("D42").Value = "'1.56"
$ws.Range("E42").Value = "  +1.08%  "

# Row 43
$ws.Range("D43").Value = "'69.67"
$ws.Range("E43").Value = "  -0.51%  "

# Row 44
$ws.Range("D44").Value = "'1.01"
$ws.Range("E44").Value = "  +0.28%  "

# Row 45
$ws.Range("D45").Value = "'0.224"
$ws.Range("E45").Value = "  -0.76%  "

# Row 46
$ws.Range("D46").Value = "'12.01"
$ws.Range("E46").Value = "  -2.55%  "

# Row 47
$ws.Range("D47").Value = "1.721.62"
$ws.Range("E47").Value = "  +8.17%  "

# Row 48
$ws.Range("D48").Value = "'109.83"
$ws.Range("E48").Value = "  -3.46%  "

# Row 49
$ws.Range("D49").Value = "'76.88"
$ws.Range("E49").Value = "  -5.83%  "

# Row 50
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.65"
$ws.Range("E50").Value = "  -3.15%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'5.14"
$ws.Range("E51").Value = "  -2.48%  "
